$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 71: "if" -> "while", 17 -> 29 ---
$ws.Range('B71').Value = "while"
$ws.Range('C71').Value = 29

# --- Pre-format the new A96:A111 cells to match the existing A-column style
#     (bold, bordered, centered) by copying the format from A95. ---
$ws.Range('A95').Copy()
$ws.Range('A96:A111').PasteSpecial(-4122)

# --- Row 96 ---
$ws.Range('A96').Value = 94
$ws.Range('B96').Value = "if"
$ws.Range('C96').Value = 17

# --- Row 97 ---
$ws.Range('A97').Value = 95
$ws.Range('B97').Value = "a"
$ws.Range('C97').Value = 81

# --- Row 98 (formula cell: "=" as a formula, matching B90's existing shape) ---
$ws.Range('A98').Value = 96
$ws.Range('B98').Formula = "=="
$ws.Range('C98').Value = 44

# --- Row 99 ---
$ws.Range('A99').Value = 97
$ws.Range('B99').Value = "a"
$ws.Range('C99').Value = 81

# --- Row 100 ---
$ws.Range('A100').Value = 98
$ws.Range('B100').Value = ":"
$ws.Range('C100').Value = 33

# --- Row 101 ---
$ws.Range('A101').Value = 99
$ws.Range('B101').Value = "if"
$ws.Range('C101').Value = 17

# --- Row 102 ---
$ws.Range('A102').Value = 100
$ws.Range('B102').Value = "a"
$ws.Range('C102').Value = 81

# --- Row 103 ---
$ws.Range('A103').Value = 101
$ws.Range('B103').Value = ">="
$ws.Range('C103').Value = 49

# --- Row 104 (literal text "5", force text so it isn't read as a number) ---
$ws.Range('A104').Value = 102
$ws.Range('B104').Value = "'5"
$ws.Range('B104').ClearFormats()
$ws.Range('C104').Value = 80

# --- Row 105 ---
$ws.Range('A105').Value = 103
$ws.Range('B105').Value = ":"
$ws.Range('C105').Value = 33

# --- Row 106 (tab + "while") ---
$ws.Range('A106').Value = 104
$ws.Range('B106').Value = "`twhile"
$ws.Range('C106').Value = -1

# --- Row 107 (literal text "True", force text so it isn't read as boolean) ---
$ws.Range('A107').Value = 105
$ws.Range('B107').Value = "'True"
$ws.Range('B107').ClearFormats()
$ws.Range('C107').Value = 83

# --- Row 108 ---
$ws.Range('A108').Value = 106
$ws.Range('B108').Value = "#"
$ws.Range('C108').Value = 43

# --- Row 109 ---
$ws.Range('A109').Value = 107
$ws.Range('B109').Value = "a"
$ws.Range('C109').Value = 81

# --- Row 110 (literal text "=", force text so it isn't read as a formula) ---
$ws.Range('A110').Value = 108
$ws.Range('B110').Value = "'="
$ws.Range('B110').ClearFormats()
$ws.Range('C110').Value = 46

# --- Row 111 (literal text "100", force text so it isn't read as a number) ---
$ws.Range('A111').Value = 109
$ws.Range('B111').Value = "'100"
$ws.Range('B111').ClearFormats()
$ws.Range('C111').Value = 80
